$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A78").Value = 77
$ws.Range("B78").Value = 1
$ws.Range("C78").Value = "2024-06-16 13:13:07"
$ws.Range("D78").Value = 200
$ws.Range("E78").Value = 6

$ws.Range("A79").Value = 78
$ws.Range("B79").Value = 2
$ws.Range("C79").Value = "2024-06-16 13:13:07"
$ws.Range("D79").Value = 200
$ws.Range("E79").Value = 0
